# Updates StructureDefinition-measure-report-evidence.xlsx:
#  - Sheet "Metadata": bump Version, Date, Publisher, replace the duplicated
#    "Contact" row with a "Jurisdiction" row, and drop the now-redundant
#    second "Contact" row.
#  - Sheet "Elements": fix the root Extension row's Short/Definition text and
#    clear the stray "N/A" mapping values on the three extension slice rows.

$wb = $excel.ActiveWorkbook

$meta = $wb.Worksheets.Item(1)
$meta.Range("B3").Value = "6.0.0"
$meta.Range("B8").Value = "2022-01-21T20:46:54+00:00"
$meta.Range("B9").Value = "Alvearie Team"
$meta.Range("A10").Value = "Jurisdiction"
$meta.Range("B10").Value = "United States of America"
$meta.Rows.Item(11).Delete()

$elements = $wb.Worksheets.Item(2)
$elements.Range("K2").Value = "Measure Report Evidence"
$elements.Range("L2").Value = "Supporting evidence showing detailed data about why a patient qualified (or did not qualify) for measure rules"
$elements.Range("AJ5").Value = ""
$elements.Range("AJ6").Value = ""
$elements.Range("AJ7").Value = ""
